$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Return Menu Music" mod row (row 26)
$ws.Range("A26").Value = "Menu Music"
$ws.Range("E26").Value = "Mirrowel"
$ws.Range("B26").Value = "return_menu_music.script"
$ws.Range("C26").Value = "scripts/"
$ws.Range("D26").Value = "Yes"

# Add the hyperlink for the mod name, same style as the other entries
$ws.Hyperlinks.Add($ws.Range("A26"), "https://www.moddb.com/mods/stalker-anomaly/addons/return-menu-music")

# Match the selection left behind by the edit
$ws.Range("E26").Select() | Out-Null
